$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.144963979721069
$ws.Range("B1").Value = 3.451892137527466
$ws.Range("C1").Value = 3.078772068023682
$ws.Range("D1").Value = 2.539528131484985
$ws.Range("E1").Value = 1.647209167480469
